$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 'Yuzhe Tang'
$ws.Range("H5").Value = 'Jesse Q. Bond'
$ws.Range("H6").Value = 'Yiyang Sun'
$ws.Range("H7").Value = 'Jeongmin Ahn'
$ws.Range("G9").Value = 'Yuzhe Tang'
$ws.Range("G10").Value = 'Yuzhe Tang'
$ws.Range("H10").Value = 'Senem Velipasalar'
$ws.Range("G11").Value = 'Jason Pollack'
$ws.Range("H14").Value = 'Era Jain'
$ws.Range("H15").Value = 'Zhenyu Gan'
$ws.Range("G16").Value = 'Yuzhe Tang'
$ws.Range("H16").Value = 'Senem Velipasalar'
$ws.Range("H18").Value = 'Yiyang Sun'
$ws.Range("G20").Value = 'Era Jain'
$ws.Range("H20").Value = 'Zhenyu Gan'
$ws.Range("G21").Value = 'Jeongmin Ahn'
$ws.Range("H21").Value = 'Svetoslava Todorova'
$ws.Range("H22").Value = 'Min Liu'
$ws.Range("G23").Value = 'Anupam Pandey'
$ws.Range("H26").Value = 'Ben Akih-Kumgeh'
$ws.Range("H27").Value = 'Jeongmin Ahn'
$ws.Range("G28").Value = 'Zhen Ma'
$ws.Range("G30").Value = 'Elizabeth Carter'
$ws.Range("H30").Value = 'Yaoying Wu'
$ws.Range("G31").Value = 'Zhenyu Gan'
$ws.Range("H32").Value = 'Era Jain'
$ws.Range("G33").Value = 'Jeongmin Ahn'
$ws.Range("H33").Value = 'Yiyang Sun'
$ws.Range("H34").Value = 'Min Liu'
$ws.Range("G35").Value = 'Nadeem Ghani'
$ws.Range("H35").Value = 'Yi Zheng'
$ws.Range("H36").Value = 'Jason Pollack'
$ws.Range("G37").Value = 'Svetoslava Todorova'
$ws.Range("G39").Value = 'Jason Pollack'
$ws.Range("H39").Value = 'Wanliang Shan'
$ws.Range("H40").Value = 'Min Liu'
$ws.Range("G41").Value = 'Zhenyu Gan'
$ws.Range("H41").Value = 'Ben Akih-Kumgeh'
$ws.Range("H43").Value = 'Ben Akih-Kumgeh'
$ws.Range("G44").Value = 'Zhen Ma'
$ws.Range("H44").Value = 'Yaoying Wu'
$ws.Range("G45").Value = 'Zhenyu Gan'
$ws.Range("H45").Value = 'Yuzhe Tang'
$ws.Range("H46").Value = 'Era Jain'
$ws.Range("G47").Value = 'Nadeem Ghani'
$ws.Range("H47").Value = 'Yi Zheng'
$ws.Range("G48").Value = 'Anupam Pandey'
$ws.Range("H48").Value = 'Yaoying Wu'
$ws.Range("G49").Value = 'Wanliang Shan'
$ws.Range("H49").Value = 'Ben Akih-Kumgeh'
$ws.Range("G50").Value = 'Senem Velipasalar'
$ws.Range("H50").Value = 'Nadeem Ghani'
$ws.Range("H51").Value = 'Yiyang Sun'
$ws.Range("H52").Value = 'Elizabeth Carter'
$ws.Range("G53").Value = 'Nadeem Ghani'
$ws.Range("H53").Value = 'Yi Zheng'
$ws.Range("H54").Value = 'Wanliang Shan'
$ws.Range("H55").Value = 'Elizabeth Carter'
$ws.Range("G56").Value = 'Zhenyu Gan'
$ws.Range("H56").Value = 'Wanliang Shan'
$ws.Range("G57").Value = 'Nadeem Ghani'
$ws.Range("H57").Value = 'Yi Zheng'
$ws.Range("G58").Value = 'Yuzhe Tang'
$ws.Range("H58").Value = 'Yiyang Sun'
$ws.Range("G59").Value = 'Ben Akih-Kumgeh'
$ws.Range("H59").Value = 'Era Jain'
$ws.Range("G60").Value = 'Nadeem Ghani'
$ws.Range("H60").Value = 'Yi Zheng'
$ws.Range("G61").Value = 'Jeongmin Ahn'
$ws.Range("H61").Value = 'Anupam Pandey'
$ws.Range("G62").Value = 'Nadeem Ghani'
$ws.Range("H62").Value = 'Yi Zheng'
$ws.Range("G64").Value = 'Jason Pollack'
$ws.Range("H64").Value = 'Nadeem Ghani'
$ws.Range("G65").Value = 'Ben Akih-Kumgeh'
$ws.Range("H65").Value = 'Ruth Chen'
$ws.Range("H66").Value = 'Yuzhe Tang'
$ws.Range("H67").Value = 'Svetoslava Todorova'
$ws.Range("H68").Value = 'Elizabeth Carter'
$ws.Range("G69").Value = 'Zhen Ma'
$ws.Range("H69").Value = 'Elizabeth Carter'
